# Updated CreateNewSupplier and DeleteASupplier tc
#
# Adds a new "DeleteASupplier" worksheet at the end of the workbook, based on
# a copy of the "CreateNewSupplier" sheet (so that it inherits the same
# column-B "best fit" width used throughout this workbook), then trims it
# down to the 2-column username/password/admin/pointofsale layout and fixes
# up the selection.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("CreateNewSupplier")

# Copy the sheet to the end of the workbook (after itself == places the
# duplicate immediately after the source sheet, i.e. at the very end).
$sourceSheet.Copy($null, $sourceSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "DeleteASupplier"

# Remove the extra companyname/firstname/lastname columns that came along
# with the copy - this sheet only needs username/password/admin/pointofsale,
# which are already present in columns A and B.
$newSheet.Range("C1:E2").ClearContents() | Out-Null

# Match the recorded selection for this sheet.
$newSheet.Range("O8").Select() | Out-Null
